# Add a new data row (row 52) to Sheet1, mirroring the row above it but
# with a new "Test Row" entry in the First Name column, and update the
# window tab-ratio / selection to reflect the edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 52

$ws.Cells.Item($row, 1).Value = 51
$ws.Cells.Item($row, 2).Value = "Test Row"
$ws.Cells.Item($row, 3).Value = "Danz"
$ws.Cells.Item($row, 4).Value = "Male"
$ws.Cells.Item($row, 5).Value = "United States"
$ws.Cells.Item($row, 6).Value = 39
$ws.Cells.Item($row, 7).Value = "15/10/2017"
$ws.Cells.Item($row, 7).NumberFormat = "MM/DD/YY"
$ws.Cells.Item($row, 8).Value = 3265

# Reflect the new selection / active cell shown in the workbook after the edit.
[void]$ws.Range("B52").Select()

# Widen the sheet-tab area in the window (cosmetic view setting).
$excel.ActiveWindow.TabRatio = 50
